$d = $word.ActiveDocument

# Remove the two trailing boilerplate paragraphs that followed the last
# "Requisitos" entry (LOT2045): the "Ver no Jupiter..." line and the
# "© 2020 ... Jekyll ..." footer line. Walk backwards so deleting a
# paragraph doesn't shift the indices of paragraphs we still need to
# inspect.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "Ver no Jupiter Salvar em pdf Salvar em docx*" -or `
        $t -like "*Powered by Jekyll and Github pages*") {
        $p.Range.Delete()
    }
}
